$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.116695
$ws.Range("H2").Value = 3.350085
$ws.Range("I2").Value = 0.008174214292497491
$ws.Range("J2").Value = 0.008174214292497492
$ws.Range("M2").Value = 0.989021
$ws.Range("N2").Value = 2.967063
$ws.Range("O2").Value = 0.05972921679266473
$ws.Range("P2").Value = 0.05972921679266473
$ws.Range("Q2").Value = 1.104434805595
$ws.Range("R2").Value = 9.939913250355
$ws.Range("S2").Value = 0.0004882394175862812
$ws.Range("T2").Value = 0.0004882394175862813
$ws.Range("G3").Value = 1.116695
$ws.Range("H3").Value = 3.350085
$ws.Range("I3").Value = 0.008174214292497491
$ws.Range("J3").Value = 0.008174214292497492
$ws.Range("O3").Value = 0.5654368392847325
$ws.Range("P3").Value = 0.5654368392847325
$ws.Range("Q3").Value = 10.45532084975166
$ws.Range("R3").Value = 94.09788764776499
$ws.Range("S3").Value = 0.004622001893185867
$ws.Range("T3").Value = 0.004622001893185868
$ws.Range("G4").Value = 1.116695
$ws.Range("H4").Value = 3.350085
$ws.Range("I4").Value = 0.008174214292497491
$ws.Range("J4").Value = 0.008174214292497492
$ws.Range("O4").Value = 0.3748339439226028
$ws.Range("P4").Value = 0.3748339439226028
$ws.Range("Q4").Value = 6.930940605225
$ws.Range("R4").Value = 62.378465447025
$ws.Range("S4").Value = 0.003063972981725342
$ws.Range("T4").Value = 0.003063972981725343
$ws.Range("I5").Value = 0.8193429796700005
$ws.Range("J5").Value = 0.8193429796700005
$ws.Range("M5").Value = 0.989021
$ws.Range("N5").Value = 2.967063
$ws.Range("O5").Value = 0.05972921679266473
$ws.Range("P5").Value = 0.05972921679266473
$ws.Range("Q5").Value = 110.7031051654733
$ws.Range("R5").Value = 996.3279464892599
$ws.Range("S5").Value = 0.04893871446025735
$ws.Range("T5").Value = 0.04893871446025735
$ws.Range("I6").Value = 0.8193429796700005
$ws.Range("J6").Value = 0.8193429796700005
$ws.Range("O6").Value = 0.5654368392847325
$ws.Range("P6").Value = 0.5654368392847325
$ws.Range("S6").Value = 0.4632867047147399
$ws.Range("T6").Value = 0.4632867047147399
$ws.Range("I7").Value = 0.8193429796700005
$ws.Range("J7").Value = 0.8193429796700005
$ws.Range("O7").Value = 0.3748339439226028
$ws.Range("P7").Value = 0.3748339439226028
$ws.Range("S7").Value = 0.3071175604950032
$ws.Range("T7").Value = 0.3071175604950032
$ws.Range("I8").Value = 0.172482806037502
$ws.Range("J8").Value = 0.1724828060375021
$ws.Range("M8").Value = 0.989021
$ws.Range("N8").Value = 2.967063
$ws.Range("O8").Value = 0.05972921679266473
$ws.Range("P8").Value = 0.05972921679266473
$ws.Range("Q8").Value = 23.304504572305
$ws.Range("R8").Value = 209.740541150745
$ws.Range("S8").Value = 0.0103022629148211
$ws.Range("T8").Value = 0.0103022629148211
$ws.Range("I9").Value = 0.172482806037502
$ws.Range("J9").Value = 0.1724828060375021
$ws.Range("O9").Value = 0.5654368392847325
$ws.Range("P9").Value = 0.5654368392847325
$ws.Range("S9").Value = 0.09752813267680674
$ws.Range("T9").Value = 0.09752813267680675
$ws.Range("I10").Value = 0.172482806037502
$ws.Range("J10").Value = 0.1724828060375021
$ws.Range("O10").Value = 0.3748339439226028
$ws.Range("P10").Value = 0.3748339439226028
$ws.Range("S10").Value = 0.06465241044587422
$ws.Range("T10").Value = 0.06465241044587422
